$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 130.90909
$ws.Range("I28").Value = 44
$ws.Range("K28").Value = 44
$ws.Range("M28").Value = 441
# Row 41
$ws.Range("H41").Value = 6428.2856
$ws.Range("I41").Value = 1333.6666
$ws.Range("J41").Value = 10249.25
$ws.Range("K41").Value = 1333.6666
$ws.Range("L41").Value = 10249.25
$ws.Range("M41").Value = -893.6666
$ws.Range("N41").Value = -11129.25
# Row 86
$ws.Range("H86").Value = 1483.6666
$ws.Range("I86").Value = 1367.3334
$ws.Range("J86").Value = 1600
$ws.Range("K86").Value = 1367.3334
$ws.Range("L86").Value = 1600
$ws.Range("M86").Value = -244.3334
$ws.Range("N86").Value = -3846
# Row 89
$ws.Range("H89").Value = 1483.6666
$ws.Range("I89").Value = 1367.3334
$ws.Range("J89").Value = 1600
$ws.Range("K89").Value = 6836.666999999999
$ws.Range("L89").Value = 8000
$ws.Range("M89").Value = -1220.666999999999
$ws.Range("N89").Value = -19232
# Row 92
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1000
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -3496
# Row 97
$ws.Range("H97").Value = 8030
$ws.Range("J97").Value = 8030
$ws.Range("L97").Value = 24090
$ws.Range("N97").Value = -25082
# Row 138
$ws.Range("H138").Value = 7129.4546
$ws.Range("J138").Value = 6352.857
$ws.Range("L138").Value = 19058.571
$ws.Range("N138").Value = -29338.571

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2855
$ws.Range("I2").Value = 2855
$ws.Range("K2").Value = 2855
$ws.Range("M2").Value = -2742
# Row 19
$ws.Range("H19").Value = 7000
$ws.Range("I19").Value = 7000
$ws.Range("K19").Value = 7000
$ws.Range("M19").Value = -6771
# Row 32
$ws.Range("H32").Value = 11310.048
$ws.Range("I32").Value = 6966.5
$ws.Range("J32").Value = 19997.143
$ws.Range("K32").Value = 6966.5
$ws.Range("L32").Value = 19997.143
$ws.Range("M32").Value = -6679.5
$ws.Range("N32").Value = -20571.143
# Row 44
$ws.Range("H44").Value = 19625.25
$ws.Range("I44").Value = 4250.5
$ws.Range("K44").Value = 4250.5
$ws.Range("M44").Value = -3762.5
# Row 88
$ws.Range("H88").Value = 3500
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3500
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3500
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4312
# Row 91
$ws.Range("H91").Value = 3500
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3500
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3500
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6308
# Row 97
$ws.Range("H97").Value = 345.17648
$ws.Range("I97").Value = 345.17648
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 345.17648
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 150.82352
$ws.Range("N97").ClearContents()
# Row 116
$ws.Range("H116").Value = 2855
$ws.Range("I116").Value = 2855
$ws.Range("K116").Value = 2855
$ws.Range("M116").Value = -561
# Row 132
$ws.Range("H132").Value = 6155.846
$ws.Range("I132").Value = 6155.846
$ws.Range("K132").Value = 18467.538
$ws.Range("M132").Value = -15937.538

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2855
$ws.Range("I3").Value = 2855
$ws.Range("K3").Value = 2855
$ws.Range("M3").Value = -2741
# Row 134
$ws.Range("H134").Value = 3501.375
$ws.Range("I134").Value = 3501.375
$ws.Range("K134").Value = 10504.125
$ws.Range("M134").Value = -7969.125

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 94
$ws.Range("H94").Value = 2000
$ws.Range("J94").Value = 2000
$ws.Range("L94").Value = 2000
$ws.Range("N94").Value = -2902
# Row 99
$ws.Range("H99").Value = 2466.5557
$ws.Range("I99").Value = 2466.5557
$ws.Range("K99").Value = 2466.5557
$ws.Range("M99").Value = -968.5556999999999
# Row 107
$ws.Range("H107").Value = 489.8125
$ws.Range("I107").Value = 364.9
$ws.Range("J107").Value = 698
$ws.Range("K107").Value = 364.9
$ws.Range("L107").Value = 698
$ws.Range("M107").Value = 1555.1
$ws.Range("N107").Value = -4538
# Row 126
$ws.Range("H126").Value = 2466.5557
$ws.Range("I126").Value = 2466.5557
$ws.Range("K126").Value = 7399.6671
$ws.Range("M126").Value = -4929.6671

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 200237.2
$ws.Range("I4").Value = 325
$ws.Range("J4").Value = 400149.4
$ws.Range("K4").Value = 975
$ws.Range("L4").Value = 1200448.2
$ws.Range("M4").Value = -863
$ws.Range("N4").Value = -1200672.2
# Row 17
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 4.5
$ws.Range("K17").Value = 13.5
$ws.Range("M17").Value = 155.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 39.23077
$ws.Range("I2").Value = 34.083332
$ws.Range("K2").Value = 34.083332
$ws.Range("M2").Value = 78.916668
# Row 22
$ws.Range("H22").Value = 2077
$ws.Range("I22").Value = 2077
$ws.Range("K22").Value = 2077
$ws.Range("M22").Value = -1548
# Row 25
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1200
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1200
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -2258
# Row 80
$ws.Range("H80").Value = 3500
$ws.Range("J80").Value = 3500
$ws.Range("L80").Value = 3500
$ws.Range("N80").Value = -5496
# Row 83
$ws.Range("H83").Value = 3500
$ws.Range("J83").Value = 3500
$ws.Range("L83").Value = 17500
$ws.Range("N83").Value = -27484
# Row 97
$ws.Range("H97").Value = 2222
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2222
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2222
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3214
# Row 132
$ws.Range("H132").Value = 3541
$ws.Range("I132").Value = 2311.5
$ws.Range("K132").Value = 6934.5
$ws.Range("M132").Value = -4404.5

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 68
$ws.Range("H68").Value = 2911.6155
$ws.Range("I68").Value = 2872.818
$ws.Range("J68").Value = 3125
$ws.Range("K68").Value = 2872.818
$ws.Range("L68").Value = 3125
$ws.Range("M68").Value = -2123.818
$ws.Range("N68").Value = -4623
# Row 71
$ws.Range("H71").Value = 2911.6155
$ws.Range("I71").Value = 2872.818
$ws.Range("J71").Value = 3125
$ws.Range("K71").Value = 14364.09
$ws.Range("L71").Value = 15625
$ws.Range("M71").Value = -10620.09
$ws.Range("N71").Value = -23113
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2872
$ws.Range("I62").Value = 2496
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 2496
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -1872
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 2872
$ws.Range("I65").Value = 2496
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 12480
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -9360
$ws.Range("N65").Value = -26240
# Row 81
$ws.Range("H81").Value = 100
$ws.Range("I81").Value = 100
$ws.Range("K81").Value = 200
$ws.Range("M81").Value = 861
# Row 84
$ws.Range("H84").Value = 100
$ws.Range("I84").Value = 100
$ws.Range("K84").Value = 1000
$ws.Range("M84").Value = 4304
# Row 113
$ws.Range("H113").Value = 673.93335
$ws.Range("I113").Value = 642.5
$ws.Range("K113").Value = 1927.5
$ws.Range("M113").Value = 242.5
# Row 122
$ws.Range("H122").Value = 224021.78
$ws.Range("I122").Value = 334199
$ws.Range("K122").Value = 1002597
$ws.Range("M122").Value = -1000147
# Row 132
$ws.Range("H132").Value = 4003.2
$ws.Range("I132").Value = 3367.2727
$ws.Range("K132").Value = 10101.8181
$ws.Range("M132").Value = -7571.8181
